$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header + text cleanup: normalize small connector words ("de", "del", "el",
# "la", "los") to title case in place names, and rename header columns to
# snake_case machine-readable names.
$cellUpdates = @(
    @{ Cell = 'A1'; Value = 'mx_state' },
    @{ Cell = 'B1'; Value = 'mx_municipality' },
    @{ Cell = 'C1'; Value = 'n_matriculas' },
    @{ Cell = 'D1'; Value = 'pct_matriculas' },
    @{ Cell = 'B4'; Value = 'Pabellón De Arteaga' },
    @{ Cell = 'B5'; Value = 'Rincón De Romos' },
    @{ Cell = 'B14'; Value = 'Chiapa De Corzo' },
    @{ Cell = 'B36'; Value = 'Hidalgo Del Parral' },
    @{ Cell = 'A48'; Value = 'Ciudad De México' },
    @{ Cell = 'A62'; Value = 'Coahuila De Zaragoza' },
    @{ Cell = 'B77'; Value = 'San Luis Del Cordero' },
    @{ Cell = 'B78'; Value = 'San Pedro Del Gallo' },
    @{ Cell = 'A82'; Value = 'Estado De México' },
    @{ Cell = 'B89'; Value = 'Ecatepec De Morelos' },
    @{ Cell = 'B92'; Value = 'Naucalpan De Juárez' },
    @{ Cell = 'B95'; Value = 'Tlalnepantla De Baz' },
    @{ Cell = 'B98'; Value = 'Apaseo El Alto' },
    @{ Cell = 'B99'; Value = 'Apaseo El Grande' },
    @{ Cell = 'B103'; Value = 'Dolores Hidalgo Cuna De La Independencia Nacional' },
    @{ Cell = 'B112'; Value = 'San Miguel De Allende' },
    @{ Cell = 'B115'; Value = 'Acapulco De Juárez' },
    @{ Cell = 'B119'; Value = 'Chilapa De Álvarez' },
    @{ Cell = 'B120'; Value = 'Chilpancingo De Los Bravo' },
    @{ Cell = 'B123'; Value = 'Coyuca De Benítez' },
    @{ Cell = 'B124'; Value = 'Iguala De La Independencia' },
    @{ Cell = 'B129'; Value = 'Técpan De Galeana' },
    @{ Cell = 'B133'; Value = 'Cuautepec De Hinojosa' },
    @{ Cell = 'B134'; Value = 'Mixquiahuala De Juárez' },
    @{ Cell = 'B135'; Value = 'Pachuca De Soto' },
    @{ Cell = 'B140'; Value = 'Tulancingo De Bravo' },
    @{ Cell = 'B143'; Value = 'Autlán De Navarro' },
    @{ Cell = 'B149'; Value = 'Lagos De Moreno' },
    @{ Cell = 'B151'; Value = 'San Diego De Alejandría' },
    @{ Cell = 'B152'; Value = 'San Juan De Los Lagos' },
    @{ Cell = 'B155'; Value = 'Tamazula De Gordiano' },
    @{ Cell = 'B158'; Value = 'Tlajomulco De Zúñiga' },
    @{ Cell = 'B160'; Value = 'Unión De Tula' },
    @{ Cell = 'A163'; Value = 'Michoacán De Ocampo' },
    @{ Cell = 'B195'; Value = 'Puente De Ixtla' },
    @{ Cell = 'B201'; Value = 'Ixtlán Del Río' },
    @{ Cell = 'B204'; Value = 'Santa María Del Oro' },
    @{ Cell = 'B209'; Value = 'San Nicolás De Los Garza' },
    @{ Cell = 'B212'; Value = 'Oaxaca De Juárez' },
    @{ Cell = 'B213'; Value = 'Ocotlán De Morelos' },
    @{ Cell = 'B218'; Value = 'San Martín De Los Cansecos' },
    @{ Cell = 'B231'; Value = 'Chalchicomula De Sesma' },
    @{ Cell = 'B244'; Value = 'Tepexi De Rodríguez' },
    @{ Cell = 'B246'; Value = 'Tlacotepec De Benito Juárez' },
    @{ Cell = 'B255'; Value = 'San Juan Del Río' },
    @{ Cell = 'B263'; Value = 'Santa María Del Río' },
    @{ Cell = 'B264'; Value = 'Villa De Ramos' },
    @{ Cell = 'B285'; Value = 'Jalpa De Méndez' },
    @{ Cell = 'B299'; Value = 'Ixtacuixtla De Mariano Matamoros' },
    @{ Cell = 'B300'; Value = 'Muñoz De Domingo Arenas' },
    @{ Cell = 'B303'; Value = 'Sanctórum De Lázaro Cárdenas' },
    @{ Cell = 'B305'; Value = 'Tepetitla De Lardizábal' },
    @{ Cell = 'A313'; Value = 'Veracruz De Ignacio De La Llave' },
    @{ Cell = 'B320'; Value = 'Hueyapan De Ocampo' },
    @{ Cell = 'B339'; Value = 'Nochistlán De Mejía' },
    @{ Cell = 'B342'; Value = 'Tlaltenango De Sánchez Román' },
    @{ Cell = 'B344'; Value = 'Villa De Cos' }
)

foreach ($update in $cellUpdates) {
    $ws.Range($update.Cell).Value = $update.Value
}

# Drop the trailing metadata/footnote rows (sample size, source, author,
# date) that used to live below the data table at rows 350-354 (row 349 was
# already blank/unused).
$ws.Rows("350:354").Delete()
